# Update the "Förändrad" (Changed) date column (C) for rows 2-15:
# change the serial date value from 45203 (2023-10-04) to 45204 (2023-10-05).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 15; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45204
    }
}
